# FIX-latitude longitude dan upload excel
# Update the template header row: column B changes from "nip" to "kelas"
# (column C already reads "email" and stays that way) and refresh the
# active cell selection to G6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "kelas"

$ws.Range("G6").Select()
